$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old rows 13-25 (will be rebuilt with new content/order)
$ws.Rows("13:25").Delete() | Out-Null

# Row 13
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Value = "01/01/2022"
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Value = "01/01/2022"
$ws.Range("A13").RowHeight = 60

# Row 14
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("A14").RowHeight = 60

# Row 15
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Value = "Programa:"
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("B15").Value = "5983729 - Fernando Vernilli Junior"
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value = "5983729 - Fernando Vernilli Junior"
$ws.Range("A15").RowHeight = 120

# Row 16
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A16").RowHeight = 120

# Row 17
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").Value = "Avaliação:"

# Row 18
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = "Método:"
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("B18").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("A18").RowHeight = 60

# Row 19
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").Value = "Critério:"
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$ws.Range("B19").Value = "a) Duas provas escritas (P1 e P2, com peso 1)b) Relatórios sobre os testes experimentais: soma das notas dos relatórios divido pelo número de relatórios (SR), com peso 1."
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$ws.Range("C19").Value = "a) Duas provas escritas (P1 e P2, com peso 1)b) Relatórios sobre os testes experimentais: soma das notas dos relatórios divido pelo número de relatórios (SR), com peso 1."
$ws.Range("A19").RowHeight = 60

# Row 20
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$ws.Range("B20").Value = "Serão aplicadas duas avaliações escritas (P1 e P2, com peso 1). A nota final serão calculada pela equaçãoNF = (P1+P2+MR)/3. NF igual ou superior a 5: aprovação direta. NF entre 3 e 4,9: recuperação. NF inferior a 3: reprovação direta."
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").Value = "Serão aplicadas duas avaliações escritas (P1 e P2, com peso 1). A nota final serão calculada pela equaçãoNF = (P1+P2+MR)/3. NF igual ou superior a 5: aprovação direta. NF entre 3 e 4,9: recuperação. NF inferior a 3: reprovação direta."
$ws.Range("A20").RowHeight = 60

# Row 21
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A21").PasteSpecial(-4122) | Out-Null
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4122) | Out-Null
$ws.Range("B21").Value = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2. Média final igual ou superior a 5 (cinco): aprovado. NF inferior a 5: reprovado."
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C21").PasteSpecial(-4122) | Out-Null
$ws.Range("C21").Value = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2. Média final igual ou superior a 5 (cinco): aprovado. NF inferior a 5: reprovado."
$ws.Range("A21").RowHeight = 120

# Row 22
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A22").PasteSpecial(-4122) | Out-Null
$ws.Range("A22").Value = "Requisitos:"

# Row 23
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
$ws.Range("B23").Value = "LOM3113 -  Tratamentos de Minérios e Hidrometalurgia  (Requisito fraco)`n"
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Value = "LOM3113 -  Tratamentos de Minérios e Hidrometalurgia  (Requisito fraco)`n"
$ws.Range("A23").RowHeight = 30

$excel.CutCopyMode = 0